# Adds a new worksheet "for_r_plot" containing DNA and PCR concentration data
# (re-shaped from Sheet1 into a "tidy" long layout for plotting in R), and
# makes it the active/selected sheet, matching the upstream commit
# "Added DNA and PCR concentrations".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Leave Sheet1's selection the way it was left when the author
#    switched away from it (selected cell G6, no longer the active tab).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("G6").Select()

# ---------------------------------------------------------------------
# 2. Add the new sheet "for_r_plot" directly after Sheet1.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "for_r_plot"

# Header row
$ws2.Range("A1").Value = "Sample_ID"
$ws2.Range("B1").Value = "Deployment"
$ws2.Range("C1").Value = "Probe"
$ws2.Range("D1").Value = "Gauze"
$ws2.Range("E1").Value = "Pooled_Tele02"
$ws2.Range("F1").Value = "Pooled_MiFishU"

# Data rows: Sample_ID, Deployment, Probe, Gauze, Pooled_Tele02, Pooled_MiFishU
# (numeric literals written without exponent notation, since the COM
# PowerShell parser here doesn't accept `E-2`-style literals)
$data = @(
    @("LE2_met_1a_01.09.23", 1, 1, 0.11, 1.58, 0.37),
    @("LE2_met_1b_01.09.23", 1, 1, 0.076799999999999993, 0.73399999999999999, 0.23400000000000001),
    @("LE2_met_1c_01.09.23", 1, 1, 0.091600000000000001, 0.94599999999999995, 0.29599999999999999),
    @("LE2_met_2a_01.09.23", 1, 2, 0.11600000000000001, 1.65, 0.41199999999999998),
    @("LE2_met_2b_01.09.23", 1, 2, 0.16200000000000001, $null, $null),
    @("LE2_met_2c_01.09.23", 1, 2, 0.073599999999999999, 2.96, 1.73),
    @("LE2_met_3a_01.09.23", 1, 3, 0.13, 1.1599999999999999, 0.32200000000000001),
    @("LE2_met_3b_01.09.23", 1, 3, 0.13, 1.61, 0.54400000000000004),
    @("LE2_met_3c_01.09.23", 1, 3, 0.129, 1.62, 0.95),
    @("LE2_met_4a_01.09.23", 2, 4, 0.098799999999999999, 3.22, 1.8),
    @("LE2_met_4b_01.09.23", 2, 4, 0.094399999999999998, 4.16, 1.18),
    @("LE2_met_4c_01.09.23", 2, 4, 0.105, 2.04, 0.33),
    @("LE2_met_5a_01.09.23", 2, 5, 0.1, 1.9, 0.22),
    @("LE2_met_5b_01.09.23", 2, 5, 0.095600000000000004, 2.68, 0.63200000000000001),
    @("LE2_met_5c_01.09.23", 2, 5, 0.083599999999999994, 1.0900000000000001, 0.17399999999999999),
    @("LE2_met_6a_01.09.23", 2, 6, 0.082799999999999999, 11.1, 0.81799999999999995),
    @("LE2_met_6b_01.09.23", 2, 6, 0.109, 1.99, 0.60599999999999998),
    @("LE2_met_6c_01.09.23", 2, 6, 0.066799999999999998, 1.7, 0.6),
    @("LE2_met_7a_02.09.23", 3, 7, 0.121, 0.27800000000000002, 0.216),
    @("LE2_met_7b_02.09.23", 3, 7, 0.071599999999999997, 0.96799999999999997, 0.41399999999999998),
    @("LE2_met_7c_02.09.23", 3, 7, 0.052400000000000002, 0.216, 0.13400000000000001),
    @("LE2_met_8a_02.09.23", 3, 8, 0.081199999999999994, 1.97, 0.312),
    @("LE2_met_8b_02.09.23", 3, 8, 0.11600000000000001, 0.79, 0.34200000000000003),
    @("LE2_met_8c_02.09.23", 3, 8, 0.56799999999999995, 0.56399999999999995, 0.64400000000000002),
    @("LE2_met_9a_02.09.23", 3, 9, 0.126, 2.02, 0.14199999999999999),
    @("LE2_met_9b_02.09.23", 3, 9, 0.155, 4.54, 0.33400000000000002),
    @("LE2_met_9c_02.09.23", 3, 9, 0.193, 3.24, 0.86),
    @("LE2_met_10a_02.09.23", 4, 10, 0.093200000000000005, 1.28, 0.39200000000000002),
    @("LE2_met_10b_02.09.23", 4, 10, 0.151, 1.23, 0.68),
    @("LE2_met_10c_02.09.23", 4, 10, 0.121, 0.48599999999999999, 0.23),
    @("LE2_met_11a_02.09.23", 4, 11, 0.16800000000000001, 1, 0.28999999999999998),
    @("LE2_met_11b_02.09.23", 4, 11, 0.106, 3.16, 0.82599999999999996),
    @("LE2_met_11c_02.09.23", 4, 11, 0.14499999999999999, 2.54, 0.34),
    @("LE2_met_12a_02.09.23", 4, 12, 0.14899999999999999, 0.59, 0.21),
    @("LE2_met_12b_02.09.23", 4, 12, 0.097600000000000006, 2.2799999999999998, 0.34599999999999997),
    @("LE2_met_12c_02.09.23", 4, 12, 0.189, 1.52, 0.222)
)

$row = 2
foreach ($rec in $data) {
    $ws2.Cells.Item($row, 1).Value = $rec[0]
    $ws2.Cells.Item($row, 2).Value = $rec[1]
    $ws2.Cells.Item($row, 3).Value = $rec[2]
    $ws2.Cells.Item($row, 4).Value = $rec[3]
    if ($null -ne $rec[4]) { $ws2.Cells.Item($row, 5).Value = $rec[4] }
    if ($null -ne $rec[5]) { $ws2.Cells.Item($row, 6).Value = $rec[5] }
    $row++
}

# Number formatting for the concentration columns (D:F), matching Sheet1's
# "0.000" custom format.
$ws2.Range("D2:F37").NumberFormat = "0.000"

# Column widths, to match the source layout.
$ws2.Columns.Item(1).ColumnWidth = 20.28515625
$ws2.Columns.Item(2).ColumnWidth = 20.28515625
$ws2.Columns.Item(3).ColumnWidth = 20.28515625
$ws2.Columns.Item(4).ColumnWidth = 26.28515625
$ws2.Columns.Item(5).ColumnWidth = 26.140625
$ws2.Columns.Item(6).ColumnWidth = 27.7109375

# Bold/shaded header row.
$ws2.Range("A1:F1").Font.Bold = $true

# ---------------------------------------------------------------------
# 3. Make the new sheet the active tab, with the cursor at C32 as left
#    by the author, topLeftCell scrolled to A3.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Application.ActiveWindow.ScrollRow = 3
$ws2.Range("C32").Select()

$wb.Save()
